$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell values to reflect the "short sign" (negative Vdelta on row 4,
# negative funding rate on row 3, etc.)
$ws.Range("D3").Value = -0.1

$ws.Range("C4").Value = -1000
$ws.Range("D4").Value = 0.2
$ws.Range("E4").Value = 16280

# Update the active selection on the sheet to H5
$ws.Range("H5").Select()
